$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Cell value updates -------------------------------------------------
# Order matters: it controls the order new entries are appended to the
# shared-strings table, which is chosen here to reproduce the target file.

$ws.Range("K4").Value = "NA"
$ws.Range("K5").Value = "NA"
$ws.Range("L4").Value = "/scratch/rtraborn/TSRchitect_plant_results/pipeline_scripts/Morton_At_PEAT_TSRchitect.Rscript "
$ws.Range("L5").Value = "/scratch/rtraborn/TSRchitect_plant_results/pipeline_scripts/Mejia-Guerra_Zm_CAGE_TSRchitect.Rscript"
$ws.Range("L2").Value = "/scratch/rtraborn/TSRchitect_plant_results/pipeline_scripts/Tokizawa_CAGE_TSRchitect.Rscript"
$ws.Range("L3").Value = "/scratch/rtraborn/TSRchitect_plant_results/pipeline_scripts/Tokizawa_Capping_TSRchitect.Rscript"

$ws.Range("K1").Value = "Name and location of Fastq Reads"

$ws.Range("J2").Value = "/scratch/rtraborn/TSRchitect_plant_results/alignment_data/Tokizawa_CAGE_align.bam"
$ws.Range("J3").Value = "/scratch/rtraborn/TSRchitect_plant_results/alignment_data/Tokizawa_Vec_Capping_align.bam"
$ws.Range("J4").Value = "/scratch/rtraborn/TSRchitect_plant_results/alignment_data/Morton_At_PEAT.bam "
$ws.Range("J5").Value = "/scratch/rtraborn/TSRchitect_plant_results/alignment_data/Zm_B73-root-1.bam, /scratch/rtraborn/TSRchitect_plant_results/alignment_data/Zm_B73-root-2.bam, /scratch/rtraborn/TSRchitect_plant_results/alignment_data/Zm_B73-shoot-1.bam,, /scratch/rtraborn/TSRchitect_plant_results/alignment_data/Zm_B73-shoot-2.bam"

$ws.Range("K3").Value = "/scratch/rtraborn/TSRchitect_plant_results/fastq_data/Tokizawa_Vec_capping_tagdust.fq"
$ws.Range("K2").Value = "/scratch/rtraborn/TSRchitect_plant_results/fastq_data/Tokizawa_CAGE_tagdust_READ1.fq; /scratch/rtraborn/TSRchitect_plant_results/fastq_data/Tokizawa_CAGE_tagdust_READ2.fq; "

# --- Column width updates ------------------------------------------------
# The engine stores width = ColumnWidth + 5/6, so subtract the padding
# back out to land as closely as possible on the target stored widths
# (78.5, 94.1640625, 87).
$ws.Columns.Item(10).ColumnWidth = 78.5 - 5/6
$ws.Columns.Item(11).ColumnWidth = 94.1640625 - 5/6
$ws.Columns.Item(12).ColumnWidth = 87 - 5/6

# --- Sheet view / selection ----------------------------------------------
$ws.Range("K5").Select()
